$d = $word.ActiveDocument

# Step 1: remove the existing _GoBack bookmark (it will be re-created later at the new location)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Step 2: update the text "... y lo muestra los datos ..." -> "... y muestra los datos ..."
$d.Content.Find.Execute(
    "El Sistema agrega un inventario y lo muestra los datos del mismo en pantalla",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El Sistema agrega un inventario y muestra los datos del mismo en pantalla",
    2
)

# Step 3: insert the _GoBack bookmark right after "... y muestra los " (before "datos del mismo...")
$r = $d.Content
$r.Find.Execute("El Sistema agrega un inventario y muestra los ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bookmarkRange = $d.Range($r.End, $r.End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
